# cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.695.71"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").Value = "1.599.46"
$ws.Range("E3").Value = "  +0.37%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("E5").Value = "  +0.28%  "

$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("E8").Value = "  +0.49%  "

$ws.Range("E9").Value = "  +1.16%  "

$ws.Range("E10").Value = "  +0.79%  "

$ws.Range("E11").Value = "  +0.97%  "

$ws.Range("D12").Value = "1.824.37"
$ws.Range("E12").Value = "  +0.38%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.04"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.73%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.571.54"
$ws.Range("E14").Value = "  -1.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.523"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.72%  "

$ws.Range("E16").Value = "  +1.42%  "

$ws.Range("D17").Value = "26.684.05"

$ws.Range("D18").Value = "0.0₃0757"
$ws.Range("E18").Value = "  +3.85%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.77"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.03%  "

$ws.Range("E20").Value = "  +0.15%  "

$ws.Range("E21").Value = "  +4.25%  "

$ws.Range("E22").Value = "  +0.94%  "

$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("E24").Value = "  +1.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.17"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.36%  "

$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.33"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.78%  "

$ws.Range("E30").Value = "  +2.66%  "

$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("E32").Value = "  +1.21%  "

$ws.Range("E33").Value = "  +1.88%  "

$ws.Range("D34").Value = "1.289.91"
$ws.Range("E34").Value = "  +0.81%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.618"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.03%  "

$ws.Range("E36").Value = "  +0.89%  "

$ws.Range("E37").Value = "  +0.70%  "

$ws.Range("E38").Value = "  +0.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.07"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +16.11%  "

$ws.Range("E40").Value = "  -1.68%  "

$ws.Range("E41").Value = "  -0.59%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.785"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.19"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.13"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.85%  "

$ws.Range("D45").Value = "1.736.83"
$ws.Range("E45").Value = "  +0.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.05"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.90%  "

$ws.Range("E47").Value = "  -0.69%  "

$ws.Range("E48").Value = "  -1.23%  "

$ws.Range("E49").Value = "  +0.75%  "

$ws.Range("E50").Value = "  +0.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.36"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.26%  "
